$d = $word.ActiveDocument

# Locate the end of the "SFX around the palace..." bullet paragraph.
$r = $d.Content
$r.Find.Execute("SFX around the palace (river, birds, footsteps, NPCs talking, fire etc)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

# Insert first new bullet: "Post processing and cool lighting"
$r.InsertParagraphAfter()
$r.Move(1, 1)
$r.InsertAfter("Post processing and cool lighting")
$r.Collapse(0)

# Insert second new bullet: "Animation of NPCs to come close to last"
$r.InsertParagraphAfter()
$r.Move(1, 1)
$r.InsertAfter("Animation of NPCs to come close to last")
